# Update the "Förändrad" (Changed) date column (C) for rows 2-23
# from 45170 (2023-09-01) to 45174 (2023-09-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
